$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VIC mystery cases")

# Update the "Date" values for the first few rows of the VIC_Mystery_cases table.
# Columns D (Exposure Date) and E (Onset of symptoms up to) are table formulas
# that recalculate automatically from column A.
$ws.Range("A2").Value = 44193
$ws.Range("A3").Value = 44191
$ws.Range("A4").Value = 44186

# Move the active selection to A5 (matches the saved selection in the file)
$ws.Range("A5").Select()
